$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.83"
$ws.Range("E2").Value = "'7.12%"
$ws.Range("D3").Value = "'40.08"
$ws.Range("E3").Value = "'7.63%"
$ws.Range("D4").Value = "'5.397"
$ws.Range("E4").Value = "'5.30%"
$ws.Range("D5").Value = "'0.08135"
$ws.Range("E5").Value = "'4.01%"
$ws.Range("D6").Value = "'4.527"
$ws.Range("E6").Value = "'2.76%"
$ws.Range("E7").Value = "'4.84%"
$ws.Range("E8").Value = "'2.41%"
$ws.Range("D10").Value = "'0.9447"
$ws.Range("E10").Value = "'2.59%"
$ws.Range("D11").Value = "'0.1359"
$ws.Range("E11").Value = "'26.11%"
$ws.Range("D12").Value = "'0.1977"
$ws.Range("E12").Value = "'4.52%"
$ws.Range("D13").Value = "'0.09287"
$ws.Range("E13").Value = "'4.53%"
$ws.Range("E14").Value = "'7.22%"
$ws.Range("D15").Value = "'0.09591"
$ws.Range("E15").Value = "'0.02%"
$ws.Range("D16").Value = "'0.001320"
$ws.Range("E16").Value = "'-4.08%"
$ws.Range("D17").Value = "'0.006034"
$ws.Range("E17").Value = "'5.69%"
$ws.Range("D18").Value = "'3.365"
$ws.Range("E18").Value = "'-0.91%"
$ws.Range("D19").Value = "'0.3525"
$ws.Range("E19").Value = "'2.91%"
$ws.Range("D20").Value = "'7.195"
$ws.Range("E20").Value = "'14.05%"
$ws.Range("E21").Value = "'2.12%"
$ws.Range("D22").Value = "'0.2450"
$ws.Range("E22").Value = "'1.42%"
$ws.Range("D23").Value = "'0.04428"
$ws.Range("E23").Value = "'1.22%"
$ws.Range("E24").Value = "'2.46%"
$ws.Range("D25").Value = "'0.004309"
$ws.Range("E25").Value = "'0.95%"
$ws.Range("E26").Value = "'-14.29%"
$ws.Range("D39").Value = "'0.02506"
$ws.Range("E39").Value = "'15.30%"
$ws.Range("D40").Value = "'0.05236"
$ws.Range("E40").Value = "'3.99%"
$ws.Range("D41").Value = "'0.007620"
$ws.Range("E41").Value = "'1.07%"
$ws.Range("D42").Value = "'0.1429"
$ws.Range("E42").Value = "'5.76%"
$ws.Range("D43").Value = "'0.009162"
$ws.Range("E43").Value = "'6.08%"
$ws.Range("E44").Value = "'4.88%"
$ws.Range("D45").Value = "'0.01077"
$ws.Range("E45").Value = "'36.34%"
$ws.Range("D46").Value = "'0.00006587"
$ws.Range("E46").Value = "'1.17%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003345"
$ws.Range("E48").Value = "'1.51%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002400"
$ws.Range("E49").Value = "'139.49%"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E51").Value = "'0.06%"
